$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('Keywords', 'Suggestions'),
  @('Fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('Order Fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('Purchase Fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('I want fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('Looking for fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('Show me fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('Fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('List of fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('Premium fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('super fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('cheap fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('Cool fans', 'We have multiple kinds of fans. Please select your choice. 1.Ceiling fans 2. Table Fans 3. Pedastal Fans 4. Wallmounted fans 5. Domestic Exhaust Fans 6. Special Fans 7. Kitchen Fans'),
  @('Ceiling fan', 'We have multiple ranges in the ceiling fans. Please select your range 1. Super Premium 2. Premium 3. Metallic/Decorative 4. Standard Plain 5. Economy 6. BEE Star rated 7. Rural 8. ISI'),
  @('Super Premium', 'Great! You are one step away. There are many models. Please select the model from the following: 1. LUSTER EROS 2. TITATNIS 3. NEBULA 4. JUPITER 5. URANUS 6. KANNON ELECTROPLATED'),
  @('Premium', 'Great! You are one step away. There are many models. Please select the model from the following: 1. KANNON PAINTED 2. TRINTON PRIME 3. AVANCER E-SENSE 4. TRINTON ELECTROPLATED 5. TRINTON PAINTED 6. STUBBY 7. RIVERA 600 MM 8. AGGNUZ PRIMUS 9. TWIRL 10. RADIANCE 11. AVANCER PRIME-ANTI DUST 12. AVANCER 13. AMOUR 14. IMPERIAL 15. LERONE 16. SENORITA 17. SENORITA 1200 18. DEVINE 19. AURA PRIME - ANTI DUST'),
  @('Metallic / Decorative', 'Great! You are one step away. There are many models. Please select the model from the following: 1. AIR 360 DECO 2. CALIBER (HI SPEED) 3. AURA 4. RIVIERA STYLE 5. FLYLEAF 6. DECAIR 7. CALIBER (HI SPEED) METALL 8. AURA METALLIC 9. HIGH SPEED DECORA METALLIC 10. HIGH SPEED DECORA 11. DIAMOND 12. FLAVIA 13. MOUNT BREESE DECO 14. HIGH SPEED METALLIC'),
  @('Standar Plain', 'Great! You are one step away. There are many models. Please select the model from the following: 1. AIR 360 DECO 2. JURA 3. HIGH BREEZE 4. HIGH SPEED 5. ULTRAA (HI SPEED) 6. MOUNT BREEZE 7. SMART BREEZE (HI SPEED) 8. RIVIERA (HI SPEED) 9. ZAPAR (HI SPEED)'),
  @('Economy', 'Great! You are one step away. There are many models. Please select the model from the following: 1. BRIZ AIR DECO 2. COOL BREEZE DECO 3. BRIZ AIR 4. COOL BREEZE 5. SEA BREEZE 6. SEA WIND 7. WINZER 1200 8. WINZER 1400 9. HILL BRIZ 10. NEO BREEZE'),
  @('BEE star rated', 'Great! You are one step away. There are many models. Please select the model from the following: 1. AURA PLUS - 5 STAR 2. HS DECORA PLUS - 5 STAR 3. HB PLUS - 5 STAR 4. COOL BREEZE DECO PLUS - 3 5. COOL BREEZE PLUS - 3 STAR 6. HS PLUS - 5 STAR 7. EMPOWER 50 - STAR 8. ENTRUST 50-5 STAR 9. ASSET 50 - 5 STAR'),
  @('Rural', 'Great! You are one step away. There are many models. Please select the model from the following: 1. PAVAN DECO 2. PAVAN'),
  @('ISI', 'Great! You are one step away. There are many models. Please select the model from the following: 1. WWHE 2. WWPS'),
  @('Table fan', 'Great! You are one step away. Please select the model from the following: 1. FIZZAIR (HIGH SPEED) 2. SDX BLACK-GOLD 3. SDX-120 4. WIND MILL 5. HIGH FLO EVA'),
  @('Pedastal fans', 'Great! You are one step away. Please select the model from the following: 1. FIZZAIR (HIGH SPEED) 2. SDX BLACK-GOLD 3. WIND MILL 4. MINI SAPPY 5. HIGH FLO EVA - WITH REM 6. HIGH FLO ESTER - WITH REM 7. HIGH FLOW LG 8. HIGH FLO WAVE 9. HIGH FLO LG 18 10. WIND FLO (HI-SPEED) 11. STORM 12. WINDSTAR FARRATA'),
  @('Wall mounted fans', 'Great! You are one step away. Please select the model from the following: 1. FIZZAIR (HIGH SPEED) 2. SDX BLACK-GOLD 3. SDX MIII 4. HIGH FLO WAVE 6. HIGH FLO LG 18 7. WIND FLO(HI-SPEED)'),
  @('Domestic exhaust fans', 'Great! Please select the model: 1. Platic Exhaust 2. Metal exhaust'),
  @('Plastic Exhaust fans', 'Great! You are one step away. Please select the model from the following: 1. MOUNT AIR 2. VENTILLUS 3. BRISK AIR HS 4. AXIAL AIR 5. BRISK AIR'),
  @('Metall exhaust', 'Great! You are one step away. Please select the model from the following: 1. TRANS AIR/ TRAMS AOR REVER 2. DRIFT AIR 3. FLUX AIR'),
  @('AIR CIRCULATOR', 'Great! You are one step away. Please select the model from the following: 1. AIR CIRCULATOR WALL 2. AIR CIRUCLATOR PEDASTAL'),
  @('Kitchen', 'Air buddy'),
  @('Show me lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('Order lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('Purchase lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('I want lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('Looking for lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('Show me lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('List of lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('Premium lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('super lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('cheap lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('Cool lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('Lamps', 'We have multiple kinds of Lightings. Please select your choice. 1. LED Lamps 2. CFL 3. FTL 4. GLS 5. Halogen 6. HID'),
  @('LED', 'Great! Please select the Wattage 1. 3W 2. 5W 3. 7W'),
  @('CFL', 'Great! You are one step away. Please select the model from the following: 1. CFL DF SPIRAL 2. CFL DF 2U COLOUR 3. CFL DF 2U 4. CFL DF T2 SPIRAL 5. CFL DF 3U 6. CFL SPIRAL G24 7. CFL DF 4U 8. CFL PIN 2U 9. CFL PIN 36W U 10. CFL-UMBRELLA-FULL SPI'),
  @('FTL', 'Great! You are one step away. Please select the model from the following: 1. SUPER SAVER 2. SUPER BRITE 3. ENERGYLUX'),
  @('GLS', 'Great! You are one step away. Please select the model from the following: 1. STANDARD GLS CLEAR LAMP 2. GLS CLEAR MUSHROOM LAMP 3. GLS MILKY LAMPS 4. SPECIAL NIGHT DECO LAMPS'),
  @('Halogen', 'Great! You are one step away. Please select the model from the following: 1. HAL 1000 2. HAL-500'),
  @('HID', 'Great! You are one step away. Please select the model from the following: 1. HPSV LAMPS 2. HPMV 3. METAL HALIDE LAMPS'),
  @('Air Cooler', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('Show me coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('Order coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('Purchase coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('I want coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('Looking for coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('Show me coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('List of coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('Premium coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('super coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('cheap coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('Cool coolers', 'We have multiple variety. Please select your choice. 1. Personal Coolers 2. Tower Coolers 3. Desert Coolers'),
  @('Personal Cooler', 'Great! You are one step away. Please select the model from the following: 1. GENIE (ACGC-PAC07) 2. MARVEL (ACGC-PAC01) 3. MARVEL DLX (ACGC-PAC201DLX)'),
  @('Tower cooler', 'Great! You are one step away. Please select the model from the following: 1. MYSTIQUE (ACGC-TAC201) 2. MYSTIQUE DLX (ACGC-TAC341) 3. TURBO 22 (ACGC-MYSTIQTRBO22) 4. TURBO 34 (ACGC-MYSTIQTRBO34) 5. TURBO 50 (ACGC-MYSTIQTRBO50)'),
  @('Desert cooler', 'Great! You are one step away. Please select the model from the following: 1. OPTIMUS 65 (ACGC-OPTIMUS65) 2. OPTIMUS 100 (ACGC-OPTIMUS100)'),
  @('Window cooler', 'Great! You are one step away. Please select the model from the following: 1. TRICOOL (ACGC-TRICOOL) 2. OZONE CHILL (ACGC-DAC541) 3. OZONE CLASSIC (ACGC-OZONECLASIC54)'),
  @('Pumps', 'We have multiple variety. Please select your choice. 1. Solar Submersible Pump 2. Solar Monobloc Pump'),
  @('Solar Submersible pump', 'Great! You are one step away. Please select the model from the following: 1. SOLAR BLDC SUBMERSIBLE PU 2. SOLAR AC SUBMERSIBLE PUMP'),
  @('Solar Monobloc Pump', 'SOLAR AC CENTRIFUGAL MONO'),
  @('Hi', 'Hi again! Happy to help! What are you looking for?'),
  @('hello', 'Hi again! Happy to help! What are you looking for?'),
  @('hey', 'Hi again! Happy to help! What are you looking for?'),
  @('Hiiiiii', 'Hi again! Happy to help! What are you looking for?'),
  @('heyyyyyy', 'Hi again! Happy to help! What are you looking for?'),
  @('oye', 'Hi again! Happy to help! What are you looking for?'),
  @('shit', 'I don''t know how to respond to that')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Reset header formatting to default (remove bold font + alignment)
$ws.Range("A1:B1").ClearFormats()

# Column widths
$ws.Columns.Item(1).ColumnWidth = 23
$ws.Columns.Item(2).ColumnWidth = 122.42578125

# Selection / view state
$ws.Activate()
$ws.Range("B79").Select() | Out-Null
